$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pi")

# Insert 9 new rows above row 122 (shifts existing 122+ down to 131+)
$ws.Rows("122:130").Insert()

# Fill in the new "install multimon-ng" section
$ws.Range("A122").Value = "install multimon-ng"
$ws.Range("B122").Value = "cd ~/temp"
$ws.Range("B122").WrapText = $true
$ws.Range("B123").Value = "git clone git://github.com/EliasOenal/multimon-ng.git"
$ws.Range("B124").Value = "cd multimon-ng"
$ws.Range("B125").Value = "mkdir build"
$ws.Range("B126").Value = "cd build"
$ws.Range("B127").Value = "qmake ../multimon-ng.pro"
$ws.Range("B128").Value = "make"
$ws.Range("B129").Value = "sudo make install"

# Update sheet view (scroll position + selection) to match post-edit state
$ws.Application.ActiveWindow.ScrollRow = 103
$ws.Range("B129").Select()

# Update workbook window size/position
$excel.ActiveWindow.Left = 19515
$excel.ActiveWindow.Top = 1365
$excel.ActiveWindow.Width = 30045
$excel.ActiveWindow.Height = 16200
